$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A4").Value = "DKS"
$ws.Range("B2").Value = "Dekstop"
$ws.Range("C2").Value = "Desktop Computer"

$ws.Range("B3").Value = "الحاسوب"
$ws.Range("C3").Value = "أجهزة الكمبيوتر المكتبية"

$ws.Range("B4").Value = "Ordinateur"
$ws.Range("C4").Value = "Ordinateurs de bureau"

$ws.Range("D10").Select() | Out-Null

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
